# Applies the "custom accuracy + 데이터 1000개" edit:
#  - Replaces the numeric data in rows 2-5 (columns A:AH) with new values
#  - Removes row 6 (dataset now has one fewer sample row)
#  - Adjusts a handful of column widths (7 <-> 8 "characters")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ColWidth($sheet, $col, $chars) {
    # Excel's ColumnWidth property is expressed in units of the workbook's
    # default font "characters"; the stored OOXML <col width> value is
    # ColumnWidth + 5/6 (the standard Normal-style column-width padding).
    # Subtracting that constant lets us land on an exact integer stored width.
    $sheet.Columns.Item($col).ColumnWidth = $chars - 5/6
}

# ---- new data values (rows 2-5, columns A..AH) ----
$row2 = @(45156.50694444445, 4.96, 4.926, 1.074, 9.69, 9.442, 3.413, 8.191000000000001, 4.526, 2.102, 4.543, 5.171, 3.55, 0.877, 3.454, 4.766, 2.485, 1.205, 0, 48.94, 9.632, 2.995, 5.545, 4.641, 0.661, 7.552, 1.905, 4.889, 2.236, 4.485, 0.29, 5.812, 2.376, 3.752)
$row3 = @(45156.51388888889, 17.026, 13.114, 0.949, 36.223, 30.713, 12.989, 45.383, 19.788, 9.037000000000001, 14.002, 14.946, 15.281, 4.129, 12.963, 18.638, 10.761, 0.731, 0.307, 192.576, 36.335, 11.864, 24.406, 13.47, 1.748, 24.434, 10.108, 10.446, 10.688, 15.567, 0.033, 40.558, 7.193, 14.874)
$row4 = @(45156.52083333334, 9.295999999999999, 7.144, 0.55, 19.626, 16.747, 7.005, 30.129, 10.696, 4.941, 7.642, 8.097, 8.308999999999999, 2.251, 6.998, 10.165, 5.858, 0.479, 0.116, 101.177, 19.818, 6.394, 13.348, 7.363, 0.945, 15.314, 5.401, 5.788, 5.768, 8.488, 0, 27.22, 3.952, 8.032)
$row5 = @(45156.52777777778, 4.46, 3.45, 0.32, 9.23, 8.029999999999999, 3.26, 16.25, 4.98, 2.36, 3.68, 3.85, 3.92, 1.07, 3.27, 4.84, 2.77, 0.34, 0.01, 44.02, 9.4, 2.97, 6.36, 3.57, 0.45, 8.17, 2.44, 2.93, 2.67, 4.06, 0, 14.65, 1.93, 3.75)

$allRows = @($row2, $row3, $row4, $row5)

for ($r = 0; $r -lt $allRows.Length; $r++) {
    $vals = $allRows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $vals[$c]
    }
}

# ---- remove row 6 (previously the last data row) ----
$ws.Rows.Item(6).Delete()

# ---- column width tweaks ----
Set-ColWidth $ws 3 8
Set-ColWidth $ws 7 8
Set-ColWidth $ws 17 8
Set-ColWidth $ws 22 8
Set-ColWidth $ws 24 7
Set-ColWidth $ws 27 8
Set-ColWidth $ws 28 8
Set-ColWidth $ws 29 8
